$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"17.73076433333334"
$ws.Range("H2").Value = [double]"53.19229300000001"
$ws.Range("I2").Value = [double]"0.004631884691211661"
$ws.Range("J2").Value = [double]"0.00463188469121166"
$ws.Range("M2").Value = [double]"281.0920463333333"
$ws.Range("N2").Value = [double]"843.2761389999999"
$ws.Range("O2").Value = [double]"0.8291026083535286"
$ws.Range("P2").Value = [double]"0.8291026083535286"
$ws.Range("Q2").Value = [double]"4983.976829510748"
$ws.Range("R2").Value = [double]"44855.79146559673"
$ws.Range("S2").Value = [double]"0.003840307679076366"
$ws.Range("T2").Value = [double]"0.003840307679076366"
# Row 3
$ws.Range("G3").Value = [double]"17.73076433333334"
$ws.Range("H3").Value = [double]"53.19229300000001"
$ws.Range("I3").Value = [double]"0.004631884691211661"
$ws.Range("J3").Value = [double]"0.00463188469121166"
$ws.Range("O3").Value = [double]"0.001324719879221983"
$ws.Range("P3").Value = [double]"0.001324719879221983"
$ws.Range("Q3").Value = [double]"7.963276338915335"
$ws.Range("R3").Value = [double]"71.66948705023802"
$ws.Range("S3").Value = [double]"6.135949728712066E-06"
$ws.Range("T3").Value = [double]"6.135949728712064E-06"
# Row 4
$ws.Range("G4").Value = [double]"17.73076433333334"
$ws.Range("H4").Value = [double]"53.19229300000001"
$ws.Range("I4").Value = [double]"0.004631884691211661"
$ws.Range("J4").Value = [double]"0.00463188469121166"
$ws.Range("M4").Value = [double]"4.452417"
$ws.Range("N4").Value = [double]"13.357251"
$ws.Range("O4").Value = [double]"0.01313274635953239"
$ws.Range("P4").Value = [double]"0.01313274635953239"
$ws.Range("Q4").Value = [double]"78.94475654072701"
$ws.Range("R4").Value = [double]"710.5028088665431"
$ws.Range("S4").Value = [double]"6.082936681628374E-05"
$ws.Range("T4").Value = [double]"6.082936681628373E-05"
# Row 5
$ws.Range("G5").Value = [double]"17.73076433333334"
$ws.Range("H5").Value = [double]"53.19229300000001"
$ws.Range("I5").Value = [double]"0.004631884691211661"
$ws.Range("J5").Value = [double]"0.00463188469121166"
$ws.Range("M5").Value = [double]"53.03808999999999"
$ws.Range("N5").Value = [double]"159.11427"
$ws.Range("O5").Value = [double]"0.156439925407717"
$ws.Range("P5").Value = [double]"0.156439925407717"
$ws.Range("Q5").Value = [double]"940.4058744801233"
$ws.Range("R5").Value = [double]"8463.65287032111"
$ws.Range("S5").Value = [double]"0.0007246116955902985"
$ws.Range("T5").Value = [double]"0.0007246116955902985"
# Row 6
$ws.Range("I6").Value = [double]"0.9353873458333681"
$ws.Range("J6").Value = [double]"0.935387345833368"
$ws.Range("M6").Value = [double]"281.0920463333333"
$ws.Range("N6").Value = [double]"843.2761389999999"
$ws.Range("O6").Value = [double]"0.8291026083535286"
$ws.Range("P6").Value = [double]"0.8291026083535286"
$ws.Range("Q6").Value = [double]"1006490.698504746"
$ws.Range("R6").Value = [double]"9058416.286542714"
$ws.Range("S6").Value = [double]"0.7755320882513296"
$ws.Range("T6").Value = [double]"0.7755320882513295"
# Row 7
$ws.Range("I7").Value = [double]"0.9353873458333681"
$ws.Range("J7").Value = [double]"0.935387345833368"
$ws.Range("O7").Value = [double]"0.001324719879221983"
$ws.Range("P7").Value = [double]"0.001324719879221983"
$ws.Range("S7").Value = [double]"0.001239126211798151"
$ws.Range("T7").Value = [double]"0.001239126211798151"
# Row 8
$ws.Range("I8").Value = [double]"0.9353873458333681"
$ws.Range("J8").Value = [double]"0.935387345833368"
$ws.Range("M8").Value = [double]"4.452417"
$ws.Range("N8").Value = [double]"13.357251"
$ws.Range("O8").Value = [double]"0.01313274635953239"
$ws.Range("P8").Value = [double]"0.01313274635953239"
$ws.Range("Q8").Value = [double]"15942.52258226556"
$ws.Range("R8").Value = [double]"143482.7032403901"
$ws.Range("S8").Value = [double]"0.01228420476074583"
$ws.Range("T8").Value = [double]"0.01228420476074583"
# Row 9
$ws.Range("I9").Value = [double]"0.9353873458333681"
$ws.Range("J9").Value = [double]"0.935387345833368"
$ws.Range("M9").Value = [double]"53.03808999999999"
$ws.Range("N9").Value = [double]"159.11427"
$ws.Range("O9").Value = [double]"0.156439925407717"
$ws.Range("P9").Value = [double]"0.156439925407717"
$ws.Range("Q9").Value = [double]"189910.5469108651"
$ws.Range("R9").Value = [double]"1709194.922197786"
$ws.Range("S9").Value = [double]"0.1463319266094945"
$ws.Range("T9").Value = [double]"0.1463319266094945"
# Row 10
$ws.Range("G10").Value = [double]"227.2177583333333"
$ws.Range("H10").Value = [double]"681.653275"
$ws.Range("I10").Value = [double]"0.0593570833501536"
$ws.Range("J10").Value = [double]"0.05935708335015359"
$ws.Range("M10").Value = [double]"281.0920463333333"
$ws.Range("N10").Value = [double]"843.2761389999999"
$ws.Range("O10").Value = [double]"0.8291026083535286"
$ws.Range("P10").Value = [double]"0.8291026083535286"
$ws.Range("Q10").Value = [double]"63869.10465318947"
$ws.Range("R10").Value = [double]"574821.9418787052"
$ws.Range("S10").Value = [double]"0.04921311262987015"
$ws.Range("T10").Value = [double]"0.04921311262987015"
# Row 11
$ws.Range("G11").Value = [double]"227.2177583333333"
$ws.Range("H11").Value = [double]"681.653275"
$ws.Range("I11").Value = [double]"0.0593570833501536"
$ws.Range("J11").Value = [double]"0.05935708335015359"
$ws.Range("O11").Value = [double]"0.001324719879221983"
$ws.Range("P11").Value = [double]"0.001324719879221983"
$ws.Range("Q11").Value = [double]"102.0484940581833"
$ws.Range("R11").Value = [double]"918.43644652365"
$ws.Range("S11").Value = [double]"7.863150828658468E-05"
$ws.Range("T11").Value = [double]"7.863150828658467E-05"
# Row 12
$ws.Range("G12").Value = [double]"227.2177583333333"
$ws.Range("H12").Value = [double]"681.653275"
$ws.Range("I12").Value = [double]"0.0593570833501536"
$ws.Range("J12").Value = [double]"0.05935708335015359"
$ws.Range("M12").Value = [double]"4.452417"
$ws.Range("N12").Value = [double]"13.357251"
$ws.Range("O12").Value = [double]"0.01313274635953239"
$ws.Range("P12").Value = [double]"0.01313274635953239"
$ws.Range("Q12").Value = [double]"1011.668209905225"
$ws.Range("R12").Value = [double]"9105.013889147025"
$ws.Range("S12").Value = [double]"0.0007795215202791901"
$ws.Range("T12").Value = [double]"0.0007795215202791902"
# Row 13
$ws.Range("G13").Value = [double]"227.2177583333333"
$ws.Range("H13").Value = [double]"681.653275"
$ws.Range("I13").Value = [double]"0.0593570833501536"
$ws.Range("J13").Value = [double]"0.05935708335015359"
$ws.Range("M13").Value = [double]"53.03808999999999"
$ws.Range("N13").Value = [double]"159.11427"
$ws.Range("O13").Value = [double]"0.156439925407717"
$ws.Range("P13").Value = [double]"0.156439925407717"
$ws.Range("Q13").Value = [double]"12051.19591608158"
$ws.Range("R13").Value = [double]"108460.7632447342"
$ws.Range("S13").Value = [double]"0.009285817691717668"
$ws.Range("T13").Value = [double]"0.009285817691717669"
# Row 14
$ws.Range("G14").Value = [double]"2.387458333333333"
$ws.Range("H14").Value = [double]"7.162374999999999"
$ws.Range("I14").Value = [double]"0.0006236861252666267"
$ws.Range("J14").Value = [double]"0.0006236861252666266"
$ws.Range("M14").Value = [double]"281.0920463333333"
$ws.Range("N14").Value = [double]"843.2761389999999"
$ws.Range("O14").Value = [double]"0.8291026083535286"
$ws.Range("P14").Value = [double]"0.8291026083535286"
$ws.Range("Q14").Value = [double]"671.095548452236"
$ws.Range("R14").Value = [double]"6039.859936070124"
$ws.Range("S14").Value = [double]"0.0005170997932524658"
$ws.Range("T14").Value = [double]"0.0005170997932524657"
# Row 15
$ws.Range("G15").Value = [double]"2.387458333333333"
$ws.Range("H15").Value = [double]"7.162374999999999"
$ws.Range("I15").Value = [double]"0.0006236861252666267"
$ws.Range("J15").Value = [double]"0.0006236861252666266"
$ws.Range("O15").Value = [double]"0.001324719879221983"
$ws.Range("P15").Value = [double]"0.001324719879221983"
$ws.Range("Q15").Value = [double]"1.072260061583333"
$ws.Range("R15").Value = [double]"9.650340554249999"
$ws.Range("S15").Value = [double]"8.262094085356326E-07"
$ws.Range("T15").Value = [double]"8.262094085356325E-07"
# Row 16
$ws.Range("G16").Value = [double]"2.387458333333333"
$ws.Range("H16").Value = [double]"7.162374999999999"
$ws.Range("I16").Value = [double]"0.0006236861252666267"
$ws.Range("J16").Value = [double]"0.0006236861252666266"
$ws.Range("M16").Value = [double]"4.452417"
$ws.Range("N16").Value = [double]"13.357251"
$ws.Range("O16").Value = [double]"0.01313274635953239"
$ws.Range("P16").Value = [double]"0.01313274635953239"
$ws.Range("Q16").Value = [double]"10.629960070125"
$ws.Range("R16").Value = [double]"95.66964063112499"
$ws.Range("S16").Value = [double]"8.190711691086153E-06"
$ws.Range("T16").Value = [double]"8.190711691086153E-06"
# Row 17
$ws.Range("G17").Value = [double]"2.387458333333333"
$ws.Range("H17").Value = [double]"7.162374999999999"
$ws.Range("I17").Value = [double]"0.0006236861252666267"
$ws.Range("J17").Value = [double]"0.0006236861252666266"
$ws.Range("M17").Value = [double]"53.03808999999999"
$ws.Range("N17").Value = [double]"159.11427"
$ws.Range("O17").Value = [double]"0.156439925407717"
$ws.Range("P17").Value = [double]"0.156439925407717"
$ws.Range("Q17").Value = [double]"71.66948705023802"
$ws.Range("R17").Value = [double]"1139.63606959125"
$ws.Range("S17").Value = [double]"9.756941091453912E-05"
$ws.Range("T17").Value = [double]"9.756941091453912E-05"
